$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in row 2
$ws.Range("B2").Value = 150
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = 1200

# Delete row 3 entirely (well #2 removed)
$ws.Rows("3:3").Delete()

# Set column B width (target stored width is 16.28515625 "characters", i.e. a
# 114px-wide column at the default Calibri 11 metrics). The host's width
# setter quantizes internally, so feed it the characters value that lands on
# the closest quantized width to the target.
$ws.Columns("B").ColumnWidth = 15.5

# Update selection to B3 (now an empty row below data)
$ws.Range("B3").Select()
